$d = $word.ActiveDocument
$xml = '<w:document xmlns:wpc="http://schemas.microsoft.com/office/word/2010/wordprocessingCanvas" xmlns:cx="http://schemas.microsoft.com/office/drawing/2014/chartex" xmlns:cx1="http://schemas.microsoft.com/office/drawing/2015/9/8/chartex" xmlns:cx2="http://schemas.microsoft.com/office/drawing/2015/10/21/chartex" xmlns:cx3="http://schemas.microsoft.com/office/drawing/2016/5/9/chartex" xmlns:cx4="http://schemas.microsoft.com/office/drawing/2016/5/10/chartex" xmlns:cx5="http://schemas.microsoft.com/office/drawing/2016/5/11/chartex" xmlns:cx6="http://schemas.microsoft.com/office/drawing/2016/5/12/chartex" xmlns:cx7="http://schemas.microsoft.com/office/drawing/2016/5/13/chartex" xmlns:cx8="http://schemas.microsoft.com/office/drawing/2016/5/14/chartex" xmlns:mc="http://schemas.openxmlformats.org/markup-compatibility/2006" xmlns:aink="http://schemas.microsoft.com/office/drawing/2016/ink" xmlns:am3d="http://schemas.microsoft.com/office/drawing/2017/model3d" xmlns:o="urn:schemas-microsoft-com:office:office" xmlns:oel="http://schemas.microsoft.com/office/2019/extlst" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math" xmlns:v="urn:schemas-microsoft-com:vml" xmlns:wp14="http://schemas.microsoft.com/office/word/2010/wordprocessingDrawing" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:w10="urn:schemas-microsoft-com:office:word" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:w15="http://schemas.microsoft.com/office/word/2012/wordml" xmlns:w16cex="http://schemas.microsoft.com/office/word/2018/wordml/cex" xmlns:w16cid="http://schemas.microsoft.com/office/word/2016/wordml/cid" xmlns:w16="http://schemas.microsoft.com/office/word/2018/wordml" xmlns:w16sdtdh="http://schemas.microsoft.com/office/word/2020/wordml/sdtdatahash" xmlns:w16se="http://schemas.microsoft.com/office/word/2015/wordml/symex" xmlns:wpg="http://schemas.microsoft.com/office/word/2010/wordprocessingGroup" xmlns:wpi="http://schemas.microsoft.com/office/word/2010/wordprocessingInk" xmlns:wne="http://schemas.microsoft.com/office/word/2006/wordml" xmlns:wps="http://schemas.microsoft.com/office/word/2010/wordprocessingShape" mc:Ignorable="w14 w15 w16se w16cid w16 w16cex w16sdtdh wp14"><w:body><w:tbl><w:tblPr><w:tblStyle w:val="LiBang"/><w:tblW w:w="0" w:type="auto"/><w:tblLook w:val="04A0" w:firstRow="1" w:lastRow="0" w:firstColumn="1" w:lastColumn="0" w:noHBand="0" w:noVBand="1"/></w:tblPr><w:tblGrid><w:gridCol w:w="1980"/><w:gridCol w:w="4253"/><w:gridCol w:w="3117"/></w:tblGrid><w:tr w:rsidR="008D50A4" w14:paraId="19CC7766" w14:textId="77777777" w:rsidTr="008E6BCE"><w:tc><w:tcPr><w:tcW w:w="1980" w:type="dxa"/><w:vAlign w:val="center"/></w:tcPr><w:p w14:paraId="387C80AC" w14:textId="77777777" w:rsidR="008D50A4" w:rsidRDefault="008D50A4" w:rsidP="008E6BCE"><w:pPr><w:jc w:val="center"/></w:pPr><w:r><w:t>Requirement ID</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="4253" w:type="dxa"/><w:vAlign w:val="center"/></w:tcPr><w:p w14:paraId="1B6117DC" w14:textId="77777777" w:rsidR="008D50A4" w:rsidRDefault="008D50A4" w:rsidP="008E6BCE"><w:pPr><w:jc w:val="center"/></w:pPr><w:r><w:t>Requirement Name</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="3117" w:type="dxa"/><w:vAlign w:val="center"/></w:tcPr><w:p w14:paraId="76DD106F" w14:textId="4B53D25D" w:rsidR="008D50A4" w:rsidRDefault="008D50A4" w:rsidP="008E6BCE"><w:pPr><w:jc w:val="center"/></w:pPr><w:r><w:t>Requirement Source</w:t></w:r></w:p></w:tc></w:tr><w:tr w:rsidR="008D50A4" w14:paraId="5D88CEE5" w14:textId="77777777" w:rsidTr="008E6BCE"><w:tc><w:tcPr><w:tcW w:w="1980" w:type="dxa"/></w:tcPr><w:p w14:paraId="540D4E4F" w14:textId="77777777" w:rsidR="008D50A4" w:rsidRDefault="008D50A4" w:rsidP="008E6BCE"><w:r w:rsidRPr="00010106"><w:t>1</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="4253" w:type="dxa"/></w:tcPr><w:p w14:paraId="1D8206C4" w14:textId="77777777" w:rsidR="008D50A4" w:rsidRDefault="008D50A4" w:rsidP="008E6BCE"><w:r w:rsidRPr="007B3683"><w:t>Sign in</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="3117" w:type="dxa"/><w:vAlign w:val="center"/></w:tcPr><w:p w14:paraId="204A4D6B" w14:textId="6935FAAC" w:rsidR="008D50A4" w:rsidRDefault="008D50A4" w:rsidP="008E6BCE"><w:r><w:t>Development Team</w:t></w:r></w:p></w:tc></w:tr><w:tr w:rsidR="008D50A4" w14:paraId="01C60C31" w14:textId="77777777" w:rsidTr="008E6BCE"><w:tc><w:tcPr><w:tcW w:w="1980" w:type="dxa"/></w:tcPr><w:p w14:paraId="64761CE7" w14:textId="77777777" w:rsidR="008D50A4" w:rsidRDefault="008D50A4" w:rsidP="008E6BCE"><w:r w:rsidRPr="00010106"><w:t>2</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="4253" w:type="dxa"/></w:tcPr><w:p w14:paraId="223FB179" w14:textId="77777777" w:rsidR="008D50A4" w:rsidRDefault="008D50A4" w:rsidP="008E6BCE"><w:r w:rsidRPr="007B3683"><w:t>Sign up</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="3117" w:type="dxa"/></w:tcPr><w:p w14:paraId="0088237C" w14:textId="56CA2764" w:rsidR="008D50A4" w:rsidRDefault="00C3327E" w:rsidP="008E6BCE"><w:r><w:t>Development Team</w:t></w:r></w:p></w:tc></w:tr><w:tr w:rsidR="008D50A4" w14:paraId="7C7D0869" w14:textId="77777777" w:rsidTr="008E6BCE"><w:tc><w:tcPr><w:tcW w:w="1980" w:type="dxa"/></w:tcPr><w:p w14:paraId="25CE5380" w14:textId="77777777" w:rsidR="008D50A4" w:rsidRDefault="008D50A4" w:rsidP="008E6BCE"><w:r w:rsidRPr="00010106"><w:t>3</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="4253" w:type="dxa"/></w:tcPr><w:p w14:paraId="77660E1A" w14:textId="77777777" w:rsidR="008D50A4" w:rsidRDefault="008D50A4" w:rsidP="008E6BCE"><w:r w:rsidRPr="007B3683"><w:t>Searching health information</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="3117" w:type="dxa"/></w:tcPr><w:p w14:paraId="213800A2" w14:textId="1676F241" w:rsidR="008D50A4" w:rsidRDefault="00C3327E" w:rsidP="008E6BCE"><w:r><w:t>Development Team</w:t></w:r></w:p></w:tc></w:tr><w:tr w:rsidR="00C3327E" w14:paraId="1A3F766E" w14:textId="77777777" w:rsidTr="008E6BCE"><w:tc><w:tcPr><w:tcW w:w="1980" w:type="dxa"/></w:tcPr><w:p w14:paraId="63B8A1DF" w14:textId="77777777" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r w:rsidRPr="00010106"><w:t>3.1</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="4253" w:type="dxa"/></w:tcPr><w:p w14:paraId="3D26E06C" w14:textId="77777777" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r w:rsidRPr="007B3683"><w:t>Search healthcare</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="3117" w:type="dxa"/></w:tcPr><w:p w14:paraId="177C651B" w14:textId="60E37573" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r><w:t>Development Team</w:t></w:r></w:p></w:tc></w:tr><w:tr w:rsidR="00C3327E" w14:paraId="3971558E" w14:textId="77777777" w:rsidTr="008E6BCE"><w:tc><w:tcPr><w:tcW w:w="1980" w:type="dxa"/></w:tcPr><w:p w14:paraId="0FAF2797" w14:textId="77777777" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r w:rsidRPr="00010106"><w:t>3.2</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="4253" w:type="dxa"/></w:tcPr><w:p w14:paraId="4D8E8639" w14:textId="77777777" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r w:rsidRPr="007B3683"><w:t>Read healthcare</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="3117" w:type="dxa"/></w:tcPr><w:p w14:paraId="6FC025A6" w14:textId="36EB9FC9" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r><w:t>Development Team</w:t></w:r></w:p></w:tc></w:tr><w:tr w:rsidR="00C3327E" w14:paraId="281CE5F6" w14:textId="77777777" w:rsidTr="008E6BCE"><w:tc><w:tcPr><w:tcW w:w="1980" w:type="dxa"/></w:tcPr><w:p w14:paraId="4BEBF2FD" w14:textId="77777777" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r w:rsidRPr="00010106"><w:t>4</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="4253" w:type="dxa"/></w:tcPr><w:p w14:paraId="35D5BBE9" w14:textId="77777777" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r w:rsidRPr="007B3683"><w:t>Managing children information</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="3117" w:type="dxa"/></w:tcPr><w:p w14:paraId="6396B7FA" w14:textId="1064086E" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r><w:t>Development Team</w:t></w:r></w:p></w:tc></w:tr><w:tr w:rsidR="00C3327E" w14:paraId="2DEEDD30" w14:textId="77777777" w:rsidTr="008E6BCE"><w:tc><w:tcPr><w:tcW w:w="1980" w:type="dxa"/></w:tcPr><w:p w14:paraId="77AED6AB" w14:textId="77777777" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r w:rsidRPr="00010106"><w:t>4.1</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="4253" w:type="dxa"/></w:tcPr><w:p w14:paraId="096E6931" w14:textId="77777777" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r w:rsidRPr="007B3683"><w:t>Search child</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="3117" w:type="dxa"/></w:tcPr><w:p w14:paraId="6BC63512" w14:textId="7D5CF03D" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r><w:t>Development Team</w:t></w:r></w:p></w:tc></w:tr><w:tr w:rsidR="00C3327E" w14:paraId="26B62528" w14:textId="77777777" w:rsidTr="008E6BCE"><w:tc><w:tcPr><w:tcW w:w="1980" w:type="dxa"/></w:tcPr><w:p w14:paraId="3FA858C6" w14:textId="77777777" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r w:rsidRPr="00010106"><w:t>4.2</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="4253" w:type="dxa"/></w:tcPr><w:p w14:paraId="386271B5" w14:textId="77777777" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r w:rsidRPr="007B3683"><w:t>Create child</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="3117" w:type="dxa"/></w:tcPr><w:p w14:paraId="33B0F510" w14:textId="578F1AFB" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r><w:t>Development Team</w:t></w:r></w:p></w:tc></w:tr><w:tr w:rsidR="00C3327E" w14:paraId="1EED879E" w14:textId="77777777" w:rsidTr="008E6BCE"><w:tc><w:tcPr><w:tcW w:w="1980" w:type="dxa"/></w:tcPr><w:p w14:paraId="34BC27CE" w14:textId="77777777" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r w:rsidRPr="00010106"><w:t>4.3</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="4253" w:type="dxa"/></w:tcPr><w:p w14:paraId="314FCD1F" w14:textId="77777777" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r w:rsidRPr="007B3683"><w:t>Read child</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="3117" w:type="dxa"/></w:tcPr><w:p w14:paraId="5D7CFC72" w14:textId="4743C9B6" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r><w:t>Development Team</w:t></w:r></w:p></w:tc></w:tr><w:tr w:rsidR="00C3327E" w14:paraId="5317E5A6" w14:textId="77777777" w:rsidTr="008E6BCE"><w:tc><w:tcPr><w:tcW w:w="1980" w:type="dxa"/></w:tcPr><w:p w14:paraId="197888F3" w14:textId="77777777" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r w:rsidRPr="00010106"><w:t>4.4</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="4253" w:type="dxa"/></w:tcPr><w:p w14:paraId="402A8FDE" w14:textId="77777777" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r w:rsidRPr="007B3683"><w:t>Update child</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="3117" w:type="dxa"/></w:tcPr><w:p w14:paraId="3BDE1D54" w14:textId="6480D6BA" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r><w:t>Development Team</w:t></w:r></w:p></w:tc></w:tr><w:tr w:rsidR="00C3327E" w14:paraId="2DDD3CF2" w14:textId="77777777" w:rsidTr="008E6BCE"><w:tc><w:tcPr><w:tcW w:w="1980" w:type="dxa"/></w:tcPr><w:p w14:paraId="56D1EBBE" w14:textId="77777777" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r w:rsidRPr="00010106"><w:t>4.5</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="4253" w:type="dxa"/></w:tcPr><w:p w14:paraId="4E4DBB14" w14:textId="77777777" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r w:rsidRPr="007B3683"><w:t>Delete child</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="3117" w:type="dxa"/></w:tcPr><w:p w14:paraId="170ECD0E" w14:textId="35C8C5B7" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r><w:t>Development Team</w:t></w:r></w:p></w:tc></w:tr><w:tr w:rsidR="00C3327E" w14:paraId="18B649D4" w14:textId="77777777" w:rsidTr="008E6BCE"><w:tc><w:tcPr><w:tcW w:w="1980" w:type="dxa"/></w:tcPr><w:p w14:paraId="1B77E5DC" w14:textId="77777777" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r w:rsidRPr="00010106"><w:t>5</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="4253" w:type="dxa"/></w:tcPr><w:p w14:paraId="61B9E982" w14:textId="77777777" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r w:rsidRPr="007B3683"><w:t>Tracking health</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="3117" w:type="dxa"/></w:tcPr><w:p w14:paraId="34401018" w14:textId="5ECE1D84" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r><w:t>Development Team</w:t></w:r></w:p></w:tc></w:tr><w:tr w:rsidR="00C3327E" w14:paraId="2639C93D" w14:textId="77777777" w:rsidTr="008E6BCE"><w:tc><w:tcPr><w:tcW w:w="1980" w:type="dxa"/></w:tcPr><w:p w14:paraId="243B0D14" w14:textId="77777777" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r w:rsidRPr="00010106"><w:t>5.1</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="4253" w:type="dxa"/></w:tcPr><w:p w14:paraId="0871818D" w14:textId="77777777" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r w:rsidRPr="007B3683"><w:t>Read tracked info</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="3117" w:type="dxa"/></w:tcPr><w:p w14:paraId="0583A411" w14:textId="0EA1E2F8" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r><w:t>Development Team</w:t></w:r></w:p></w:tc></w:tr><w:tr w:rsidR="00C3327E" w14:paraId="70BD5ABF" w14:textId="77777777" w:rsidTr="008E6BCE"><w:tc><w:tcPr><w:tcW w:w="1980" w:type="dxa"/></w:tcPr><w:p w14:paraId="11E9F152" w14:textId="77777777" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r w:rsidRPr="00010106"><w:t>5.2</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="4253" w:type="dxa"/></w:tcPr><w:p w14:paraId="56DB8731" w14:textId="77777777" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r w:rsidRPr="007B3683"><w:t>Update tracked info</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="3117" w:type="dxa"/></w:tcPr><w:p w14:paraId="3EDE49AD" w14:textId="040D962D" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r><w:t>Development Team</w:t></w:r></w:p></w:tc></w:tr><w:tr w:rsidR="00C3327E" w14:paraId="44E26C72" w14:textId="77777777" w:rsidTr="008E6BCE"><w:tc><w:tcPr><w:tcW w:w="1980" w:type="dxa"/></w:tcPr><w:p w14:paraId="2B4F369E" w14:textId="77777777" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r w:rsidRPr="00010106"><w:t>5.3</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="4253" w:type="dxa"/></w:tcPr><w:p w14:paraId="18EF0E2A" w14:textId="77777777" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r w:rsidRPr="007B3683"><w:t>Search child for tracked</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="3117" w:type="dxa"/></w:tcPr><w:p w14:paraId="25FF0EAB" w14:textId="76EF265A" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r><w:t>Development Team</w:t></w:r></w:p></w:tc></w:tr><w:tr w:rsidR="00C3327E" w14:paraId="49290BA1" w14:textId="77777777" w:rsidTr="008E6BCE"><w:tc><w:tcPr><w:tcW w:w="1980" w:type="dxa"/></w:tcPr><w:p w14:paraId="3D671F57" w14:textId="77777777" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r w:rsidRPr="00010106"><w:t>6</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="4253" w:type="dxa"/></w:tcPr><w:p w14:paraId="382A22D9" w14:textId="77777777" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r w:rsidRPr="007B3683"><w:t>Storing health records</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="3117" w:type="dxa"/></w:tcPr><w:p w14:paraId="4387BB4B" w14:textId="2C9BF038" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r><w:t>Development Team</w:t></w:r></w:p></w:tc></w:tr><w:tr w:rsidR="00C3327E" w14:paraId="248361F7" w14:textId="77777777" w:rsidTr="008E6BCE"><w:tc><w:tcPr><w:tcW w:w="1980" w:type="dxa"/></w:tcPr><w:p w14:paraId="46B00218" w14:textId="77777777" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r w:rsidRPr="00010106"><w:t>6.1</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="4253" w:type="dxa"/></w:tcPr><w:p w14:paraId="694D2E0E" w14:textId="77777777" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r w:rsidRPr="007B3683"><w:t>Medical history</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="3117" w:type="dxa"/></w:tcPr><w:p w14:paraId="50ECB909" w14:textId="4B51E447" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r><w:t>Development Team</w:t></w:r></w:p></w:tc></w:tr><w:tr w:rsidR="00C3327E" w14:paraId="607A9573" w14:textId="77777777" w:rsidTr="008E6BCE"><w:tc><w:tcPr><w:tcW w:w="1980" w:type="dxa"/></w:tcPr><w:p w14:paraId="5F521672" w14:textId="77777777" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r w:rsidRPr="00010106"><w:t>6.1.1</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="4253" w:type="dxa"/></w:tcPr><w:p w14:paraId="6B340BD2" w14:textId="77777777" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r w:rsidRPr="007B3683"><w:t>Search history</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="3117" w:type="dxa"/></w:tcPr><w:p w14:paraId="46B67385" w14:textId="249EFCBB" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r><w:t>Development Team</w:t></w:r></w:p></w:tc></w:tr><w:tr w:rsidR="00C3327E" w14:paraId="19C2F15A" w14:textId="77777777" w:rsidTr="008E6BCE"><w:tc><w:tcPr><w:tcW w:w="1980" w:type="dxa"/></w:tcPr><w:p w14:paraId="59C43857" w14:textId="77777777" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r w:rsidRPr="00010106"><w:t>6.1.2</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="4253" w:type="dxa"/></w:tcPr><w:p w14:paraId="57B23A18" w14:textId="77777777" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r w:rsidRPr="007B3683"><w:t>Create history</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="3117" w:type="dxa"/></w:tcPr><w:p w14:paraId="6B6B9869" w14:textId="2D2BB18E" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r><w:t>Development Team</w:t></w:r></w:p></w:tc></w:tr><w:tr w:rsidR="00C3327E" w14:paraId="229843CE" w14:textId="77777777" w:rsidTr="008E6BCE"><w:tc><w:tcPr><w:tcW w:w="1980" w:type="dxa"/></w:tcPr><w:p w14:paraId="66B1D140" w14:textId="77777777" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r w:rsidRPr="00010106"><w:t>6.1.3</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="4253" w:type="dxa"/></w:tcPr><w:p w14:paraId="4EB56FFF" w14:textId="77777777" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r w:rsidRPr="007B3683"><w:t>Read history</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="3117" w:type="dxa"/></w:tcPr><w:p w14:paraId="21966B1F" w14:textId="4BE35C9F" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r><w:t>Development Team</w:t></w:r></w:p></w:tc></w:tr><w:tr w:rsidR="00C3327E" w14:paraId="564E4D69" w14:textId="77777777" w:rsidTr="008E6BCE"><w:tc><w:tcPr><w:tcW w:w="1980" w:type="dxa"/></w:tcPr><w:p w14:paraId="466DEA83" w14:textId="77777777" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r w:rsidRPr="00010106"><w:t>6.1.4</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="4253" w:type="dxa"/></w:tcPr><w:p w14:paraId="01CA71AB" w14:textId="77777777" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r w:rsidRPr="007B3683"><w:t>Update history</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="3117" w:type="dxa"/></w:tcPr><w:p w14:paraId="4E776AF0" w14:textId="7F1DD748" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r><w:t>Development Team</w:t></w:r></w:p></w:tc></w:tr><w:tr w:rsidR="00C3327E" w14:paraId="2972FDFD" w14:textId="77777777" w:rsidTr="008E6BCE"><w:tc><w:tcPr><w:tcW w:w="1980" w:type="dxa"/></w:tcPr><w:p w14:paraId="14FC85EC" w14:textId="77777777" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r w:rsidRPr="00010106"><w:t>6.1.5</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="4253" w:type="dxa"/></w:tcPr><w:p w14:paraId="36AA5ADE" w14:textId="77777777" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r w:rsidRPr="007B3683"><w:t>Delete history</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="3117" w:type="dxa"/></w:tcPr><w:p w14:paraId="406CCFF7" w14:textId="3ABE4257" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r><w:t>Development Team</w:t></w:r></w:p></w:tc></w:tr><w:tr w:rsidR="00C3327E" w14:paraId="625EF061" w14:textId="77777777" w:rsidTr="008E6BCE"><w:tc><w:tcPr><w:tcW w:w="1980" w:type="dxa"/></w:tcPr><w:p w14:paraId="3B86A2C9" w14:textId="77777777" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r w:rsidRPr="00010106"><w:t>6.1.6</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="4253" w:type="dxa"/></w:tcPr><w:p w14:paraId="371FC4A0" w14:textId="77777777" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r w:rsidRPr="007B3683"><w:t>Search child for history</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="3117" w:type="dxa"/></w:tcPr><w:p w14:paraId="58C6EC2A" w14:textId="1F5BCFB3" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r><w:t>Development Team</w:t></w:r></w:p></w:tc></w:tr><w:tr w:rsidR="00C3327E" w14:paraId="6C29A254" w14:textId="77777777" w:rsidTr="008E6BCE"><w:tc><w:tcPr><w:tcW w:w="1980" w:type="dxa"/></w:tcPr><w:p w14:paraId="6F38B35E" w14:textId="77777777" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r w:rsidRPr="00010106"><w:t>6.2</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="4253" w:type="dxa"/></w:tcPr><w:p w14:paraId="6D462032" w14:textId="77777777" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r w:rsidRPr="007B3683"><w:t>Test results</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="3117" w:type="dxa"/></w:tcPr><w:p w14:paraId="364E1287" w14:textId="321C8F70" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r><w:t>Development Team</w:t></w:r></w:p></w:tc></w:tr><w:tr w:rsidR="00C3327E" w14:paraId="5B62A62E" w14:textId="77777777" w:rsidTr="008E6BCE"><w:tc><w:tcPr><w:tcW w:w="1980" w:type="dxa"/></w:tcPr><w:p w14:paraId="14454FD4" w14:textId="77777777" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r w:rsidRPr="00010106"><w:t>6.2.1</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="4253" w:type="dxa"/></w:tcPr><w:p w14:paraId="42142D4A" w14:textId="77777777" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r w:rsidRPr="007B3683"><w:t>Search result</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="3117" w:type="dxa"/></w:tcPr><w:p w14:paraId="1B1A4F58" w14:textId="763E7843" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r><w:t>Development Team</w:t></w:r></w:p></w:tc></w:tr><w:tr w:rsidR="00C3327E" w14:paraId="0D8C8B22" w14:textId="77777777" w:rsidTr="008E6BCE"><w:tc><w:tcPr><w:tcW w:w="1980" w:type="dxa"/></w:tcPr><w:p w14:paraId="728E09ED" w14:textId="77777777" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r w:rsidRPr="00010106"><w:t>6.2.2</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="4253" w:type="dxa"/></w:tcPr><w:p w14:paraId="4E1FB36C" w14:textId="77777777" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r w:rsidRPr="007B3683"><w:t>Create result</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="3117" w:type="dxa"/></w:tcPr><w:p w14:paraId="3133067C" w14:textId="57C68C87" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r><w:t>Development Team</w:t></w:r></w:p></w:tc></w:tr><w:tr w:rsidR="00C3327E" w14:paraId="76385691" w14:textId="77777777" w:rsidTr="008E6BCE"><w:tc><w:tcPr><w:tcW w:w="1980" w:type="dxa"/></w:tcPr><w:p w14:paraId="63C240C9" w14:textId="77777777" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r w:rsidRPr="00010106"><w:t>6.2.3</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="4253" w:type="dxa"/></w:tcPr><w:p w14:paraId="2E3C4A05" w14:textId="77777777" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r w:rsidRPr="007B3683"><w:t>Read result</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="3117" w:type="dxa"/></w:tcPr><w:p w14:paraId="59F063FA" w14:textId="52768AC4" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r><w:t>Development Team</w:t></w:r></w:p></w:tc></w:tr><w:tr w:rsidR="00C3327E" w14:paraId="0009A2FF" w14:textId="77777777" w:rsidTr="008E6BCE"><w:tc><w:tcPr><w:tcW w:w="1980" w:type="dxa"/></w:tcPr><w:p w14:paraId="6A5DF654" w14:textId="77777777" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r w:rsidRPr="00010106"><w:lastRenderedPageBreak/><w:t>6.2.4</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="4253" w:type="dxa"/></w:tcPr><w:p w14:paraId="31BFE5E2" w14:textId="77777777" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r w:rsidRPr="007B3683"><w:t>Update result</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="3117" w:type="dxa"/></w:tcPr><w:p w14:paraId="29E2E220" w14:textId="7781BEA6" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r><w:t>Development Team</w:t></w:r></w:p></w:tc></w:tr><w:tr w:rsidR="00C3327E" w14:paraId="11FEB011" w14:textId="77777777" w:rsidTr="008E6BCE"><w:tc><w:tcPr><w:tcW w:w="1980" w:type="dxa"/></w:tcPr><w:p w14:paraId="0CBEBB0D" w14:textId="77777777" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r w:rsidRPr="00010106"><w:t>6.2.5</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="4253" w:type="dxa"/></w:tcPr><w:p w14:paraId="0C1A8DDC" w14:textId="77777777" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r w:rsidRPr="007B3683"><w:t>Delete result</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="3117" w:type="dxa"/></w:tcPr><w:p w14:paraId="71C0F187" w14:textId="736A0775" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r><w:t>Development Team</w:t></w:r></w:p></w:tc></w:tr><w:tr w:rsidR="00C3327E" w14:paraId="2A7E98AF" w14:textId="77777777" w:rsidTr="008E6BCE"><w:tc><w:tcPr><w:tcW w:w="1980" w:type="dxa"/></w:tcPr><w:p w14:paraId="72B52C5A" w14:textId="77777777" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r w:rsidRPr="00010106"><w:t>6.2.6</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="4253" w:type="dxa"/></w:tcPr><w:p w14:paraId="2B737108" w14:textId="77777777" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r w:rsidRPr="007B3683"><w:t>Search child for result</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="3117" w:type="dxa"/></w:tcPr><w:p w14:paraId="06C40EF6" w14:textId="1C8EF06E" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r><w:t>Development Team</w:t></w:r></w:p></w:tc></w:tr><w:tr w:rsidR="00C3327E" w14:paraId="691C6C60" w14:textId="77777777" w:rsidTr="008E6BCE"><w:tc><w:tcPr><w:tcW w:w="1980" w:type="dxa"/></w:tcPr><w:p w14:paraId="55B43695" w14:textId="77777777" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r w:rsidRPr="00010106"><w:t>6.3</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="4253" w:type="dxa"/></w:tcPr><w:p w14:paraId="7D9ED721" w14:textId="77777777" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r w:rsidRPr="007B3683"><w:t>Health issues</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="3117" w:type="dxa"/></w:tcPr><w:p w14:paraId="51B0077F" w14:textId="507796B1" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r><w:t>Development Team</w:t></w:r></w:p></w:tc></w:tr><w:tr w:rsidR="00C3327E" w14:paraId="379B38DA" w14:textId="77777777" w:rsidTr="008E6BCE"><w:tc><w:tcPr><w:tcW w:w="1980" w:type="dxa"/></w:tcPr><w:p w14:paraId="02768588" w14:textId="77777777" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r w:rsidRPr="00010106"><w:t>6.3.1</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="4253" w:type="dxa"/></w:tcPr><w:p w14:paraId="17ED600A" w14:textId="77777777" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r w:rsidRPr="007B3683"><w:t>Search issue</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="3117" w:type="dxa"/></w:tcPr><w:p w14:paraId="088A2ADE" w14:textId="3182456A" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r><w:t>Development Team</w:t></w:r></w:p></w:tc></w:tr><w:tr w:rsidR="00C3327E" w14:paraId="5DC3558C" w14:textId="77777777" w:rsidTr="008E6BCE"><w:tc><w:tcPr><w:tcW w:w="1980" w:type="dxa"/></w:tcPr><w:p w14:paraId="58E0D17F" w14:textId="77777777" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r w:rsidRPr="00010106"><w:t>6.3.2</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="4253" w:type="dxa"/></w:tcPr><w:p w14:paraId="43D1C7DC" w14:textId="77777777" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r w:rsidRPr="007B3683"><w:t>Create issue</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="3117" w:type="dxa"/></w:tcPr><w:p w14:paraId="61ED0AC7" w14:textId="40BCB4B8" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r><w:t>Development Team</w:t></w:r></w:p></w:tc></w:tr><w:tr w:rsidR="00C3327E" w14:paraId="56AA0FFC" w14:textId="77777777" w:rsidTr="008E6BCE"><w:tc><w:tcPr><w:tcW w:w="1980" w:type="dxa"/></w:tcPr><w:p w14:paraId="29B5D339" w14:textId="77777777" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r w:rsidRPr="00010106"><w:t>6.3.3</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="4253" w:type="dxa"/></w:tcPr><w:p w14:paraId="3F74F2BF" w14:textId="77777777" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r w:rsidRPr="007B3683"><w:t>Read issue</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="3117" w:type="dxa"/></w:tcPr><w:p w14:paraId="04915658" w14:textId="68AF6978" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r><w:t>Development Team</w:t></w:r></w:p></w:tc></w:tr><w:tr w:rsidR="00C3327E" w14:paraId="638A75CA" w14:textId="77777777" w:rsidTr="008E6BCE"><w:tc><w:tcPr><w:tcW w:w="1980" w:type="dxa"/></w:tcPr><w:p w14:paraId="55342238" w14:textId="77777777" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r w:rsidRPr="00010106"><w:t>6.3.4</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="4253" w:type="dxa"/></w:tcPr><w:p w14:paraId="4C1198DB" w14:textId="77777777" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r w:rsidRPr="007B3683"><w:t>Update issue</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="3117" w:type="dxa"/></w:tcPr><w:p w14:paraId="5A21F395" w14:textId="16CD90B7" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r><w:t>Development Team</w:t></w:r></w:p></w:tc></w:tr><w:tr w:rsidR="00C3327E" w14:paraId="3D56D63A" w14:textId="77777777" w:rsidTr="008E6BCE"><w:tc><w:tcPr><w:tcW w:w="1980" w:type="dxa"/></w:tcPr><w:p w14:paraId="51D5BE58" w14:textId="77777777" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r w:rsidRPr="00010106"><w:t>6.3.5</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="4253" w:type="dxa"/></w:tcPr><w:p w14:paraId="0FF40865" w14:textId="77777777" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r w:rsidRPr="007B3683"><w:t>Delete issue</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="3117" w:type="dxa"/></w:tcPr><w:p w14:paraId="1FF2B6AD" w14:textId="0C47EB50" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r><w:t>Development Team</w:t></w:r></w:p></w:tc></w:tr><w:tr w:rsidR="00C3327E" w14:paraId="1905742C" w14:textId="77777777" w:rsidTr="008E6BCE"><w:tc><w:tcPr><w:tcW w:w="1980" w:type="dxa"/></w:tcPr><w:p w14:paraId="47653E92" w14:textId="77777777" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r w:rsidRPr="00010106"><w:t>6.3.6</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="4253" w:type="dxa"/></w:tcPr><w:p w14:paraId="121E362B" w14:textId="77777777" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r w:rsidRPr="007B3683"><w:t>Search child for issue</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="3117" w:type="dxa"/></w:tcPr><w:p w14:paraId="2C10BBE8" w14:textId="2E1BFF39" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r><w:t>Development Team</w:t></w:r></w:p></w:tc></w:tr><w:tr w:rsidR="00C3327E" w14:paraId="45CBE3A2" w14:textId="77777777" w:rsidTr="008E6BCE"><w:tc><w:tcPr><w:tcW w:w="1980" w:type="dxa"/></w:tcPr><w:p w14:paraId="01467B01" w14:textId="77777777" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r w:rsidRPr="00010106"><w:t>7</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="4253" w:type="dxa"/></w:tcPr><w:p w14:paraId="07E97FA0" w14:textId="77777777" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r w:rsidRPr="007B3683"><w:t>Consulting healthcare</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="3117" w:type="dxa"/></w:tcPr><w:p w14:paraId="23D1E1B0" w14:textId="3AF845D5" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r><w:t>Development Team</w:t></w:r></w:p></w:tc></w:tr><w:tr w:rsidR="00C3327E" w14:paraId="2A1F1A55" w14:textId="77777777" w:rsidTr="008E6BCE"><w:tc><w:tcPr><w:tcW w:w="1980" w:type="dxa"/></w:tcPr><w:p w14:paraId="39E89302" w14:textId="77777777" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r w:rsidRPr="00010106"><w:t>7.1</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="4253" w:type="dxa"/></w:tcPr><w:p w14:paraId="2A5C8B18" w14:textId="77777777" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r w:rsidRPr="007B3683"><w:t>Consultation history</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="3117" w:type="dxa"/></w:tcPr><w:p w14:paraId="3B50E99E" w14:textId="2F45AFC0" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r><w:t>Development Team</w:t></w:r></w:p></w:tc></w:tr><w:tr w:rsidR="00C3327E" w14:paraId="3C8AE64F" w14:textId="77777777" w:rsidTr="008E6BCE"><w:tc><w:tcPr><w:tcW w:w="1980" w:type="dxa"/></w:tcPr><w:p w14:paraId="0641F0BF" w14:textId="77777777" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r w:rsidRPr="00010106"><w:t>7.2</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="4253" w:type="dxa"/></w:tcPr><w:p w14:paraId="2ABD75E7" w14:textId="77777777" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r w:rsidRPr="007B3683"><w:t>Create consultation</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="3117" w:type="dxa"/></w:tcPr><w:p w14:paraId="4160E055" w14:textId="43842D8B" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r><w:t>Development Team</w:t></w:r></w:p></w:tc></w:tr><w:tr w:rsidR="00C3327E" w14:paraId="3E73851B" w14:textId="77777777" w:rsidTr="008E6BCE"><w:tc><w:tcPr><w:tcW w:w="1980" w:type="dxa"/></w:tcPr><w:p w14:paraId="37B23ED7" w14:textId="77777777" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r w:rsidRPr="00010106"><w:t>8</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="4253" w:type="dxa"/></w:tcPr><w:p w14:paraId="4C218549" w14:textId="77777777" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r w:rsidRPr="007B3683"><w:t>Connecting with professionals</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="3117" w:type="dxa"/></w:tcPr><w:p w14:paraId="3D0330CA" w14:textId="1F7273B4" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r><w:t>Development Team</w:t></w:r></w:p></w:tc></w:tr><w:tr w:rsidR="00C3327E" w14:paraId="78205828" w14:textId="77777777" w:rsidTr="008E6BCE"><w:tc><w:tcPr><w:tcW w:w="1980" w:type="dxa"/></w:tcPr><w:p w14:paraId="70C640D6" w14:textId="77777777" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r w:rsidRPr="00010106"><w:t>8.1</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="4253" w:type="dxa"/></w:tcPr><w:p w14:paraId="00B23E64" w14:textId="77777777" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r w:rsidRPr="007B3683"><w:t>Connection history</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="3117" w:type="dxa"/></w:tcPr><w:p w14:paraId="5110447D" w14:textId="5B52A141" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r><w:t>Development Team</w:t></w:r></w:p></w:tc></w:tr><w:tr w:rsidR="00C3327E" w14:paraId="7459092E" w14:textId="77777777" w:rsidTr="008E6BCE"><w:tc><w:tcPr><w:tcW w:w="1980" w:type="dxa"/></w:tcPr><w:p w14:paraId="0BE2FFA7" w14:textId="77777777" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r w:rsidRPr="00010106"><w:t>8.2</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="4253" w:type="dxa"/></w:tcPr><w:p w14:paraId="1E766CBB" w14:textId="77777777" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r w:rsidRPr="007B3683"><w:t>Create connection</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="3117" w:type="dxa"/></w:tcPr><w:p w14:paraId="40360B37" w14:textId="245FEAF0" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r><w:t>Development Team</w:t></w:r></w:p></w:tc></w:tr><w:tr w:rsidR="00C3327E" w14:paraId="322550FD" w14:textId="77777777" w:rsidTr="008E6BCE"><w:tc><w:tcPr><w:tcW w:w="1980" w:type="dxa"/></w:tcPr><w:p w14:paraId="16C21F80" w14:textId="77777777" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r w:rsidRPr="00010106"><w:t>9</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="4253" w:type="dxa"/></w:tcPr><w:p w14:paraId="303F39A8" w14:textId="77777777" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r w:rsidRPr="007B3683"><w:t>Reminding schedule</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="3117" w:type="dxa"/></w:tcPr><w:p w14:paraId="0493B8CD" w14:textId="2E39F83A" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r><w:t>Development Team</w:t></w:r></w:p></w:tc></w:tr><w:tr w:rsidR="00C3327E" w14:paraId="5049EF4A" w14:textId="77777777" w:rsidTr="008E6BCE"><w:tc><w:tcPr><w:tcW w:w="1980" w:type="dxa"/></w:tcPr><w:p w14:paraId="5CA0BC72" w14:textId="77777777" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r w:rsidRPr="00010106"><w:t>9.1</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="4253" w:type="dxa"/></w:tcPr><w:p w14:paraId="1D6A46A5" w14:textId="77777777" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r w:rsidRPr="007B3683"><w:t>Reminder Switch</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="3117" w:type="dxa"/></w:tcPr><w:p w14:paraId="1124EAA1" w14:textId="42D5C4EE" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r><w:t>Development Team</w:t></w:r></w:p></w:tc></w:tr><w:tr w:rsidR="00C3327E" w14:paraId="046EB939" w14:textId="77777777" w:rsidTr="008E6BCE"><w:tc><w:tcPr><w:tcW w:w="1980" w:type="dxa"/></w:tcPr><w:p w14:paraId="3F6DCAB1" w14:textId="77777777" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r w:rsidRPr="00010106"><w:t>9.2</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="4253" w:type="dxa"/></w:tcPr><w:p w14:paraId="7205B5D9" w14:textId="77777777" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r w:rsidRPr="007B3683"><w:t>Manager health check-up</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="3117" w:type="dxa"/></w:tcPr><w:p w14:paraId="601E45F5" w14:textId="379616B4" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r><w:t>Development Team</w:t></w:r></w:p></w:tc></w:tr><w:tr w:rsidR="00C3327E" w14:paraId="3CCAD32B" w14:textId="77777777" w:rsidTr="008E6BCE"><w:tc><w:tcPr><w:tcW w:w="1980" w:type="dxa"/></w:tcPr><w:p w14:paraId="26C484CB" w14:textId="77777777" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r w:rsidRPr="00010106"><w:t>9.2.1</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="4253" w:type="dxa"/></w:tcPr><w:p w14:paraId="5C27255A" w14:textId="77777777" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r w:rsidRPr="007B3683"><w:t>Search check-up</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="3117" w:type="dxa"/></w:tcPr><w:p w14:paraId="64D0E415" w14:textId="3DC9A38E" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r><w:t>Development Team</w:t></w:r></w:p></w:tc></w:tr><w:tr w:rsidR="00C3327E" w14:paraId="2D7CDE31" w14:textId="77777777" w:rsidTr="008E6BCE"><w:tc><w:tcPr><w:tcW w:w="1980" w:type="dxa"/></w:tcPr><w:p w14:paraId="4C13F4C8" w14:textId="77777777" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r w:rsidRPr="00010106"><w:t>9.2.2</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="4253" w:type="dxa"/></w:tcPr><w:p w14:paraId="34A28688" w14:textId="77777777" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r w:rsidRPr="007B3683"><w:t>Create check-up</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="3117" w:type="dxa"/></w:tcPr><w:p w14:paraId="5A0FDC22" w14:textId="05FCD4AB" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r><w:t>Development Team</w:t></w:r></w:p></w:tc></w:tr><w:tr w:rsidR="00C3327E" w14:paraId="774B6BB3" w14:textId="77777777" w:rsidTr="008E6BCE"><w:tc><w:tcPr><w:tcW w:w="1980" w:type="dxa"/></w:tcPr><w:p w14:paraId="2A64CF79" w14:textId="77777777" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r w:rsidRPr="00010106"><w:t>9.2.3</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="4253" w:type="dxa"/></w:tcPr><w:p w14:paraId="78782E58" w14:textId="77777777" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r w:rsidRPr="007B3683"><w:t>Read check-up</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="3117" w:type="dxa"/></w:tcPr><w:p w14:paraId="3596DA8F" w14:textId="3170CE5B" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r><w:t>Development Team</w:t></w:r></w:p></w:tc></w:tr><w:tr w:rsidR="00C3327E" w14:paraId="6FA1AD05" w14:textId="77777777" w:rsidTr="008E6BCE"><w:tc><w:tcPr><w:tcW w:w="1980" w:type="dxa"/></w:tcPr><w:p w14:paraId="236FA45F" w14:textId="77777777" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r w:rsidRPr="00010106"><w:t>9.2.4</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="4253" w:type="dxa"/></w:tcPr><w:p w14:paraId="0D018B2F" w14:textId="77777777" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r w:rsidRPr="007B3683"><w:t>Update check-up</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="3117" w:type="dxa"/></w:tcPr><w:p w14:paraId="0F45462C" w14:textId="2E794C3F" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r><w:t>Development Team</w:t></w:r></w:p></w:tc></w:tr><w:tr w:rsidR="00C3327E" w14:paraId="026F8AC9" w14:textId="77777777" w:rsidTr="008E6BCE"><w:tc><w:tcPr><w:tcW w:w="1980" w:type="dxa"/></w:tcPr><w:p w14:paraId="28778B02" w14:textId="77777777" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r w:rsidRPr="00010106"><w:t>9.2.5</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="4253" w:type="dxa"/></w:tcPr><w:p w14:paraId="0D1C9EAD" w14:textId="77777777" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r w:rsidRPr="007B3683"><w:t>Delete check-up</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="3117" w:type="dxa"/></w:tcPr><w:p w14:paraId="1E00093F" w14:textId="7AEA750A" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r><w:t>Development Team</w:t></w:r></w:p></w:tc></w:tr><w:tr w:rsidR="00C3327E" w14:paraId="12771183" w14:textId="77777777" w:rsidTr="008E6BCE"><w:tc><w:tcPr><w:tcW w:w="1980" w:type="dxa"/></w:tcPr><w:p w14:paraId="0C2BECC8" w14:textId="77777777" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r w:rsidRPr="00010106"><w:t>9.2.6</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="4253" w:type="dxa"/></w:tcPr><w:p w14:paraId="0C276B23" w14:textId="77777777" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r w:rsidRPr="007B3683"><w:t>Search child for check-up</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="3117" w:type="dxa"/></w:tcPr><w:p w14:paraId="5F45C3D8" w14:textId="79EA78B2" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r><w:t>Development Team</w:t></w:r></w:p></w:tc></w:tr><w:tr w:rsidR="00C3327E" w14:paraId="4731B78E" w14:textId="77777777" w:rsidTr="008E6BCE"><w:tc><w:tcPr><w:tcW w:w="1980" w:type="dxa"/></w:tcPr><w:p w14:paraId="7CE78B94" w14:textId="77777777" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r w:rsidRPr="00010106"><w:t>9.3</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="4253" w:type="dxa"/></w:tcPr><w:p w14:paraId="59E2A76C" w14:textId="77777777" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r w:rsidRPr="007B3683"><w:t>Manager vaccination</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="3117" w:type="dxa"/></w:tcPr><w:p w14:paraId="673DC169" w14:textId="36967891" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r><w:t>Development Team</w:t></w:r></w:p></w:tc></w:tr><w:tr w:rsidR="00C3327E" w14:paraId="50BE23F1" w14:textId="77777777" w:rsidTr="008E6BCE"><w:tc><w:tcPr><w:tcW w:w="1980" w:type="dxa"/></w:tcPr><w:p w14:paraId="595B9754" w14:textId="77777777" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r w:rsidRPr="00010106"><w:t>9.3.1</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="4253" w:type="dxa"/></w:tcPr><w:p w14:paraId="29C95C11" w14:textId="77777777" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r w:rsidRPr="007B3683"><w:t>Search vaccination</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="3117" w:type="dxa"/></w:tcPr><w:p w14:paraId="2CDDB6E4" w14:textId="13BDB19E" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r><w:t>Development Team</w:t></w:r></w:p></w:tc></w:tr><w:tr w:rsidR="00C3327E" w14:paraId="2F0FB694" w14:textId="77777777" w:rsidTr="008E6BCE"><w:tc><w:tcPr><w:tcW w:w="1980" w:type="dxa"/></w:tcPr><w:p w14:paraId="2975A999" w14:textId="77777777" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r w:rsidRPr="00010106"><w:t xml:space="preserve">9.3.2 </w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="4253" w:type="dxa"/></w:tcPr><w:p w14:paraId="21F80C47" w14:textId="77777777" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r w:rsidRPr="007B3683"><w:t>Create vaccination</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="3117" w:type="dxa"/></w:tcPr><w:p w14:paraId="43FBD463" w14:textId="7C0AF076" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r><w:t>Development Team</w:t></w:r></w:p></w:tc></w:tr><w:tr w:rsidR="00C3327E" w14:paraId="037D1C9C" w14:textId="77777777" w:rsidTr="008E6BCE"><w:tc><w:tcPr><w:tcW w:w="1980" w:type="dxa"/></w:tcPr><w:p w14:paraId="08DC12BA" w14:textId="77777777" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r w:rsidRPr="00010106"><w:lastRenderedPageBreak/><w:t>9.3.3</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="4253" w:type="dxa"/></w:tcPr><w:p w14:paraId="7DC71D91" w14:textId="77777777" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r w:rsidRPr="007B3683"><w:t>Read vaccination</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="3117" w:type="dxa"/></w:tcPr><w:p w14:paraId="5AC64620" w14:textId="09BF1542" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r><w:t>Development Team</w:t></w:r></w:p></w:tc></w:tr><w:tr w:rsidR="00C3327E" w14:paraId="474F96BD" w14:textId="77777777" w:rsidTr="008E6BCE"><w:tc><w:tcPr><w:tcW w:w="1980" w:type="dxa"/></w:tcPr><w:p w14:paraId="71BF3E1F" w14:textId="77777777" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r w:rsidRPr="00010106"><w:t>9.3.4</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="4253" w:type="dxa"/></w:tcPr><w:p w14:paraId="50BF26F8" w14:textId="77777777" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r w:rsidRPr="007B3683"><w:t>Update vaccination</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="3117" w:type="dxa"/></w:tcPr><w:p w14:paraId="5046CDBC" w14:textId="2FF53AF2" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r><w:t>Development Team</w:t></w:r></w:p></w:tc></w:tr><w:tr w:rsidR="00C3327E" w14:paraId="3A2B3082" w14:textId="77777777" w:rsidTr="008E6BCE"><w:tc><w:tcPr><w:tcW w:w="1980" w:type="dxa"/></w:tcPr><w:p w14:paraId="5356AAA6" w14:textId="77777777" w:rsidR="00C3327E" w:rsidRPr="0075169B" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r w:rsidRPr="00010106"><w:t>9.3.5</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="4253" w:type="dxa"/></w:tcPr><w:p w14:paraId="3EDD81EA" w14:textId="77777777" w:rsidR="00C3327E" w:rsidRPr="00DA396F" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r w:rsidRPr="007B3683"><w:t>Delete vaccination</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="3117" w:type="dxa"/></w:tcPr><w:p w14:paraId="52CDBF28" w14:textId="3D1F6525" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r><w:t>Development Team</w:t></w:r></w:p></w:tc></w:tr><w:tr w:rsidR="00C3327E" w14:paraId="44C903C1" w14:textId="77777777" w:rsidTr="008E6BCE"><w:tc><w:tcPr><w:tcW w:w="1980" w:type="dxa"/></w:tcPr><w:p w14:paraId="25F0AFF0" w14:textId="77777777" w:rsidR="00C3327E" w:rsidRPr="0075169B" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r w:rsidRPr="00010106"><w:t>9.3.6</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="4253" w:type="dxa"/></w:tcPr><w:p w14:paraId="1D9A5352" w14:textId="77777777" w:rsidR="00C3327E" w:rsidRPr="00DA396F" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r w:rsidRPr="007B3683"><w:t>Search child for vaccination</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="3117" w:type="dxa"/></w:tcPr><w:p w14:paraId="447B2F0E" w14:textId="32CA9085" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r><w:t>Development Team</w:t></w:r></w:p></w:tc></w:tr><w:tr w:rsidR="00C3327E" w14:paraId="55BC1A7C" w14:textId="77777777" w:rsidTr="008E6BCE"><w:tc><w:tcPr><w:tcW w:w="1980" w:type="dxa"/></w:tcPr><w:p w14:paraId="6EA1CA84" w14:textId="77777777" w:rsidR="00C3327E" w:rsidRPr="0075169B" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r w:rsidRPr="00010106"><w:t xml:space="preserve">10 </w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="4253" w:type="dxa"/></w:tcPr><w:p w14:paraId="4FE32E20" w14:textId="77777777" w:rsidR="00C3327E" w:rsidRPr="00DA396F" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r w:rsidRPr="007B3683"><w:t>Viewing health posts</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="3117" w:type="dxa"/></w:tcPr><w:p w14:paraId="596D6B59" w14:textId="2CC23959" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r><w:t>Development Team</w:t></w:r></w:p></w:tc></w:tr><w:tr w:rsidR="00C3327E" w14:paraId="4D566E56" w14:textId="77777777" w:rsidTr="008E6BCE"><w:tc><w:tcPr><w:tcW w:w="1980" w:type="dxa"/></w:tcPr><w:p w14:paraId="2C50E921" w14:textId="77777777" w:rsidR="00C3327E" w:rsidRPr="0075169B" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r w:rsidRPr="00010106"><w:t>10.1</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="4253" w:type="dxa"/></w:tcPr><w:p w14:paraId="370C67E9" w14:textId="77777777" w:rsidR="00C3327E" w:rsidRPr="00DA396F" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r w:rsidRPr="007B3683"><w:t>Search post</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="3117" w:type="dxa"/></w:tcPr><w:p w14:paraId="71FB3613" w14:textId="6530E754" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r><w:t>Development Team</w:t></w:r></w:p></w:tc></w:tr><w:tr w:rsidR="00C3327E" w14:paraId="05CD9AFA" w14:textId="77777777" w:rsidTr="008E6BCE"><w:tc><w:tcPr><w:tcW w:w="1980" w:type="dxa"/></w:tcPr><w:p w14:paraId="07D3C7C4" w14:textId="77777777" w:rsidR="00C3327E" w:rsidRPr="0075169B" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r w:rsidRPr="00010106"><w:t>10.2</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="4253" w:type="dxa"/></w:tcPr><w:p w14:paraId="62A35A9B" w14:textId="77777777" w:rsidR="00C3327E" w:rsidRPr="00DA396F" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r w:rsidRPr="007B3683"><w:t>Read post</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="3117" w:type="dxa"/></w:tcPr><w:p w14:paraId="2FD7D330" w14:textId="23DF6197" w:rsidR="00C3327E" w:rsidRDefault="00C3327E" w:rsidP="00C3327E"><w:r><w:t>Development Team</w:t></w:r></w:p></w:tc></w:tr></w:tbl><w:p w14:paraId="62AC944A" w14:textId="77777777" w:rsidR="00F36DA5" w:rsidRDefault="00F36DA5"/><w:sectPr w:rsidR="00F36DA5"><w:pgSz w:w="12240" w:h="15840"/><w:pgMar w:top="1440" w:right="1440" w:bottom="1440" w:left="1440" w:header="720" w:footer="720" w:gutter="0"/><w:cols w:space="720"/><w:docGrid w:linePitch="360"/></w:sectPr></w:body></w:document>'
$d.Content.InsertXML($xml)
